# Regenerate merged AHB files
#  - rename the "_old" / "_new" header suffixes to the specific version
#    tags ("_FV2310" / "_FV2404")
#  - freeze the header row (row 1)
#  - wrap the used range in an Excel Table ("Table1") with its AutoFilter

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header row renames: columns A:J were "<Name>_old", columns L:U were
#    "<Name>_new" (column K is the literal "diff" header and is untouched).
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"

$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# 2) Freeze the header row (pane split below row 1).
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3) Turn the used range into a proper Table (ListObject) with its AutoFilter.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U55"), 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

Write-Host "done"
